$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "Turn on Draft Countdown Timer" ---
# Goes in right after the existing "Schedule 8 AM Tuesday Email..." row (row 3),
# pushing the Thurs/Fri Auto-Draft-Log rows (and everything below) down by one.
$ws.Rows(4).Insert()
$ws.Range("A4").Value = "Turn on Draft Countdown Timer"
$ws.Range("B4").Value = 45468
$ws.Range("C4").Value = 0.375
$ws.Range("D4").Value = 45468
$ws.Range("E4").Value = 0.41666666666666669
$ws.Range("F4").Value = $false

# --- Insert a new row for "Turn off Draft Countdown Timer" ---
# Goes in right after the Thurs/Fri rows (now rows 5 & 6), before the
# RTTBC Post Draft Tasks block, pushing that block down by one more.
$ws.Rows(7).Insert()
$ws.Range("A7").Value = "Turn off Draft Countdown Timer"
$ws.Range("B7").Value = 45472
$ws.Range("C7").Value = 0.45833333333333331
$ws.Range("D7").Value = 45472
$ws.Range("E7").Value = 0.5
$ws.Range("F7").Value = $false

# --- Adjust the time slot for "RTTBC Post Draft Data Entry Validation" ---
# (now on row 9 after the two inserts above) so it no longer overlaps the
# new countdown-timer entries.
$ws.Range("C9").Value = 0.45833333333333331
$ws.Range("E9").Value = 0.5

# --- Restore the user's on-sheet selection ---
$ws.Range("A15").Select()
